# Insert a new "snapshot" timestamp column before the trailing "nom" /
# "url_produit" columns (previously BX/BY, now shifted to BY/BZ).
#
# Before: ... BW(last price timestamp) | BX=nom | BY=url_produit
# After:  ... BW(last price timestamp) | BX=new timestamp | BY=nom | BZ=url_produit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column 76 is "BX" - inserting here shifts the existing BX/BY (nom/url_produit)
# one column to the right, to BY/BZ, and keeps their styles/values intact.
$ws.Columns.Item(76).Insert()

# New header cell for the inserted column.
$ws.Cells.Item(1, 76).Value2 = "2026-01-31 04:56:21"

# For the data rows that already have a running price history (rows 2-80),
# the new snapshot simply repeats the most recent known price, taken from
# column BW (75), which is unaffected by the insert above.
for ($row = 2; $row -le 80; $row++) {
    $price = $ws.Cells.Item($row, 75).Value2
    $ws.Cells.Item($row, 76).Value2 = $price
}

# Rows 81-206 have no price history yet, so the newly inserted cell in
# column BX (76) is left blank for them, matching the rest of that row.
